$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append the new sentence to the paragraph ending in
#    "Vi startede på 3 og sidste sprint i dette projekt."
# ------------------------------------------------------------------
$r = $d.Content
[void]$r.Find.Execute("Vi startede på 3 og sidste sprint i dette projekt.", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $r.Paragraphs(1)
$targetRange = $targetPara.Range
[void]$targetRange.MoveEnd(1, -1)   # exclude the paragraph mark
[void]$targetRange.Collapse(0)      # collapse to the end of the visible text
[void]$targetRange.InsertAfter(" Rapporten blev sat sammen så den bare manglede korrektur og at blive sat pænt op. Kodemæssigt blev sammenligningen gjort færdig og der mangler diverse udregninger, oprydning osv.")

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark so it ends up at the end of the
#    "12-12-2017: Tolvte dag" paragraph (right after "dag").
#    (A collapsed range sitting exactly at end-of-paragraph-text
#    confuses Bookmarks.Add in this host, so a harmless temp marker
#    is appended/removed to keep the insertion point mid-text while
#    the bookmark is created.)
# ------------------------------------------------------------------
$r2 = $d.Content
[void]$r2.Find.Execute("12-12-2017", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$datePara = $r2.Paragraphs(1)
$dateRange = $datePara.Range
[void]$dateRange.MoveEnd(1, -1)     # exclude the paragraph mark
$endPos = $dateRange.End

$marker = $d.Range($endPos, $endPos)
[void]$marker.InsertAfter("ZZMARKZZ")

$bmPoint = $d.Range($endPos, $endPos)
[void]$d.Bookmarks.Add("_GoBack", $bmPoint)

$markerRange = $d.Range($endPos, $endPos + 8)
$markerRange.Text = ""

# ------------------------------------------------------------------
# 3) Insert a brand-new paragraph right after "12-12-2017: Tolvte dag"
#    containing Christian's remark.
# ------------------------------------------------------------------
$r3 = $d.Content
[void]$r3.Find.Execute("12-12-2017", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$datePara2 = $r3.Paragraphs(1)
$dateRange2 = $datePara2.Range
[void]$dateRange2.InsertParagraphAfter()

$r4 = $d.Content
[void]$r4.Find.Execute("12-12-2017", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$datePara3 = $r4.Paragraphs(1)
$newPara = $datePara3.Next()
[void]$newPara.Range.InsertAfter("Christian var meget positiv overrasket over hele gruppens resultat. Intet at sætte en finger på.")
